$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -19.46322899669338
$ws.Range("C2").Value = -0.8794961225143775
$ws.Range("D2").Value = -19.46322899669338
$ws.Range("E2").Value = -19.46322899669338
$ws.Range("F2").Value = -19.46322899669338
$ws.Range("G2").Value = -19.46322899669338
$ws.Range("H2").Value = -19.46322899669338
$ws.Range("I2").Value = -19.46322899669338
$ws.Range("J2").Value = -19.46322899669338
$ws.Range("K2").Value = -19.46322899669338

$ws.Range("B3").Value = -19.46322899669338
$ws.Range("C3").Value = -19.46322899669338
$ws.Range("D3").Value = -19.46322899669338
$ws.Range("E3").Value = -19.46322899669338
$ws.Range("F3").Value = -19.46322899669338
$ws.Range("G3").Value = -19.46322899669338
$ws.Range("H3").Value = -19.46322899669338
$ws.Range("I3").Value = 0.370237122757831
$ws.Range("J3").Value = -19.46322899669338
$ws.Range("K3").Value = -19.46322899669338

$ws.Range("B4").Value = -19.46322899669338
$ws.Range("C4").Value = -0.8498387856437518
$ws.Range("D4").Value = -0.2619545791724876
$ws.Range("E4").Value = -19.46322899669338
$ws.Range("F4").Value = 4.321926198693706
$ws.Range("G4").Value = -19.46322899669338
$ws.Range("H4").Value = 2.055224865740906
$ws.Range("I4").Value = -19.46322899669338
$ws.Range("J4").Value = 2.810041021169531
$ws.Range("K4").Value = -19.46322899669338

$ws.Range("B5").Value = -19.46322899669338
$ws.Range("C5").Value = -0.04516022229034487
$ws.Range("D5").Value = -19.46322899669338
$ws.Range("E5").Value = -19.46322899669338
$ws.Range("F5").Value = -19.46322899669338
$ws.Range("G5").Value = 3.574044633765279
$ws.Range("H5").Value = -19.46322899669338
$ws.Range("I5").Value = -19.46322899669338
$ws.Range("J5").Value = -19.46322899669338
$ws.Range("K5").Value = -19.46322899669338

$ws.Range("B6").Value = -19.46322899669338
$ws.Range("C6").Value = -19.46322899669338
$ws.Range("D6").Value = -19.46322899669338
$ws.Range("E6").Value = -19.46322899669338
$ws.Range("F6").Value = -19.46322899669338
$ws.Range("G6").Value = -19.46322899669338
$ws.Range("H6").Value = -19.46322899669338
$ws.Range("I6").Value = -19.46322899669338
$ws.Range("J6").Value = -19.46322899669338
$ws.Range("K6").Value = -19.46322899669338

$ws.Range("B7").Value = 3.257488917230996
$ws.Range("C7").Value = -19.46322899669338
$ws.Range("D7").Value = -19.46322899669338
$ws.Range("E7").Value = -19.46322899669338
$ws.Range("F7").Value = -19.46322899669338
$ws.Range("G7").Value = -19.46322899669338
$ws.Range("H7").Value = -19.46322899669338
$ws.Range("I7").Value = -19.46322899669338
$ws.Range("J7").Value = -19.46322899669338
$ws.Range("K7").Value = -19.46322899669338

$ws.Range("B8").Value = -19.46322899669338
$ws.Range("C8").Value = -19.46322899669338
$ws.Range("D8").Value = -19.46322899669338
$ws.Range("E8").Value = 1.432125137147397
$ws.Range("F8").Value = -19.46322899669338
$ws.Range("G8").Value = -19.46322899669338
$ws.Range("H8").Value = -19.46322899669338
$ws.Range("I8").Value = -19.46322899669338
$ws.Range("J8").Value = -19.46322899669338
$ws.Range("K8").Value = -19.46322899669338

$ws.Range("B9").Value = 3.383608242619341
$ws.Range("C9").Value = -19.46322899669338
$ws.Range("D9").Value = -19.46322899669338
$ws.Range("E9").Value = -19.46322899669338
$ws.Range("F9").Value = -19.46322899669338
$ws.Range("G9").Value = -19.46322899669338
$ws.Range("H9").Value = -19.46322899669338
$ws.Range("I9").Value = -19.46322899669338
$ws.Range("J9").Value = -19.46322899669338
$ws.Range("K9").Value = -19.46322899669338

$ws.Range("B10").Value = -19.46322899669338
$ws.Range("C10").Value = -19.46322899669338
$ws.Range("D10").Value = -19.46322899669338
$ws.Range("E10").Value = -19.46322899669338
$ws.Range("F10").Value = -19.46322899669338
$ws.Range("G10").Value = -19.46322899669338
$ws.Range("H10").Value = -19.46322899669338
$ws.Range("I10").Value = 0.3666191389117253
$ws.Range("J10").Value = -19.46322899669338
$ws.Range("K10").Value = 2.004323321347989

$ws.Range("B11").Value = -19.46322899669338
$ws.Range("C11").Value = -19.46322899669338
$ws.Range("D11").Value = -19.46322899669338
$ws.Range("E11").Value = 2.386733822027332
$ws.Range("F11").Value = -19.46322899669338
$ws.Range("G11").Value = 1.385886662630315
$ws.Range("H11").Value = -19.46322899669338
$ws.Range("I11").Value = -19.46322899669338
$ws.Range("J11").Value = -19.46322899669338
$ws.Range("K11").Value = 1.353977187739261

$ws.Range("B12").Value = -19.46322899669338
$ws.Range("C12").Value = -19.46322899669338
$ws.Range("D12").Value = -19.46322899669338
$ws.Range("E12").Value = -19.46322899669338
$ws.Range("F12").Value = -19.46322899669338
$ws.Range("G12").Value = -19.46322899669338
$ws.Range("H12").Value = -19.46322899669338
$ws.Range("I12").Value = -19.46322899669338
$ws.Range("J12").Value = -19.46322899669338
$ws.Range("K12").Value = -19.46322899669338

$ws.Range("B13").Value = -19.46322899669338
$ws.Range("C13").Value = -19.46322899669338
$ws.Range("D13").Value = -19.46322899669338
$ws.Range("E13").Value = 2.065055790301249
$ws.Range("F13").Value = -19.46322899669338
$ws.Range("G13").Value = -19.46322899669338
$ws.Range("H13").Value = -19.46322899669338
$ws.Range("I13").Value = -19.46322899669338
$ws.Range("J13").Value = 0.8852020084444313
$ws.Range("K13").Value = 2.741073778357798

$ws.Range("B14").Value = -19.46322899669338
$ws.Range("C14").Value = -19.46322899669338
$ws.Range("D14").Value = 1.030834545610962
$ws.Range("E14").Value = -19.46322899669338
$ws.Range("F14").Value = -19.46322899669338
$ws.Range("G14").Value = -19.46322899669338
$ws.Range("H14").Value = -19.46322899669338
$ws.Range("I14").Value = -19.46322899669338
$ws.Range("J14").Value = -19.46322899669338
$ws.Range("K14").Value = 1.717627740900721

$ws.Range("B15").Value = -19.46322899669338
$ws.Range("C15").Value = -19.46322899669338
$ws.Range("D15").Value = -0.5665861695257596
$ws.Range("E15").Value = -19.46322899669338
$ws.Range("F15").Value = -19.46322899669338
$ws.Range("G15").Value = -19.46322899669338
$ws.Range("H15").Value = -19.46322899669338
$ws.Range("I15").Value = -19.46322899669338
$ws.Range("J15").Value = -19.46322899669338
$ws.Range("K15").Value = -19.46322899669338

$ws.Range("B16").Value = -19.46322899669338
$ws.Range("C16").Value = -19.46322899669338
$ws.Range("D16").Value = -19.46322899669338
$ws.Range("E16").Value = -19.46322899669338
$ws.Range("F16").Value = -19.46322899669338
$ws.Range("G16").Value = -19.46322899669338
$ws.Range("H16").Value = -19.46322899669338
$ws.Range("I16").Value = -19.46322899669338
$ws.Range("J16").Value = 2.113497524597812
$ws.Range("K16").Value = -19.46322899669338

$ws.Range("B17").Value = -19.46322899669338
$ws.Range("C17").Value = 0.2773439944810916
$ws.Range("D17").Value = -0.5494834780758469
$ws.Range("E17").Value = -19.46322899669338
$ws.Range("F17").Value = -19.46322899669338
$ws.Range("G17").Value = -19.46322899669338
$ws.Range("H17").Value = 2.131438673829223
$ws.Range("I17").Value = 0.1842874000456073
$ws.Range("J17").Value = 1.761110133568005
$ws.Range("K17").Value = -19.46322899669338

$ws.Range("B18").Value = -19.46322899669338
$ws.Range("C18").Value = -19.46322899669338
$ws.Range("D18").Value = -19.46322899669338
$ws.Range("E18").Value = -19.46322899669338
$ws.Range("F18").Value = -19.46322899669338
$ws.Range("G18").Value = -19.46322899669338
$ws.Range("H18").Value = 2.112371645063319
$ws.Range("I18").Value = -0.5098936108915334
$ws.Range("J18").Value = 1.775232306730592
$ws.Range("K18").Value = -19.46322899669338

$ws.Range("B19").Value = -19.46322899669338
$ws.Range("C19").Value = -19.46322899669338
$ws.Range("D19").Value = 2.96032628839318
$ws.Range("E19").Value = -19.46322899669338
$ws.Range("F19").Value = -19.46322899669338
$ws.Range("G19").Value = -19.46322899669338
$ws.Range("H19").Value = 1.728004359108158
$ws.Range("I19").Value = 1.443089861171301
$ws.Range("J19").Value = -19.46322899669338
$ws.Range("K19").Value = -19.46322899669338

$ws.Range("B20").Value = -19.46322899669338
$ws.Range("C20").Value = 3.356020996911492
$ws.Range("D20").Value = 2.996631076542754
$ws.Range("E20").Value = -19.46322899669338
$ws.Range("F20").Value = -19.46322899669338
$ws.Range("G20").Value = -19.46322899669338
$ws.Range("H20").Value = 0.9285484633203589
$ws.Range("I20").Value = 3.684878404994943
$ws.Range("J20").Value = -19.46322899669338
$ws.Range("K20").Value = 1.789598679013684

$ws.Range("B21").Value = -19.46322899669338
$ws.Range("C21").Value = 2.696284225357511
$ws.Range("D21").Value = -19.46322899669338
$ws.Range("E21").Value = 2.979534582223521
$ws.Range("F21").Value = -19.46322899669338
$ws.Range("G21").Value = 2.453417697811908
$ws.Range("H21").Value = 0.9428233804342481
$ws.Range("I21").Value = -19.46322899669338
$ws.Range("J21").Value = -19.46322899669338
$ws.Range("K21").Value = -19.46322899669338
